$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Index" column header (A1) to "i" (shared string + table1 column name
# update automatically).
$ws.Range("A1").Value = "i"

# Re-index the "testdata" table's index column from 0 instead of 1: every data row
# A2:A503 decreases by one (row 2 -> 0, row 3 -> 1, ... row 503 -> 501).
$lastRow = 503
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Narrow column A now that the longest value lost a digit.
$ws.Columns.Item(1).ColumnWidth = 3.14
